# Adds the missing xG / goals columns (D:G) for the six Verona matches that
# were still lacking that data (rows 10-15 of the sheet), bringing the sheet
# up to date with the latest scraped results.
#
# Row  Home                Away        xG_home    xG_away    goals_home goals_away
# 10   Atalanta            Verona      1.20139    1.78911    0          2
# 11   Verona              Cagliari    3.02671    1.67924    1          1
# 12   Lazio               Verona      0.647828   0.955041   1          2
# 13   Verona              Sampdoria   1.49333    0.779474   1          2
# 14   Fiorentina          Verona      1.68253    1.24682    1          1
# 15   Verona              Inter       1.23629    2.0312     1          2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @{ Row = 10; D = "1.20139";  E = "1.78911";  F = "0"; G = "2" },
    @{ Row = 11; D = "3.02671";  E = "1.67924";  F = "1"; G = "1" },
    @{ Row = 12; D = "0.647828"; E = "0.955041"; F = "1"; G = "2" },
    @{ Row = 13; D = "1.49333";  E = "0.779474"; F = "1"; G = "2" },
    @{ Row = 14; D = "1.68253";  E = "1.24682";  F = "1"; G = "1" },
    @{ Row = 15; D = "1.23629";  E = "2.0312";   F = "1"; G = "2" }
)

foreach ($entry in $newData) {
    $r = $entry.Row
    foreach ($col in @("D", "E", "F", "G")) {
        $cell = $ws.Range("$col$r")
        # Prefix with a literal apostrophe so the numeric-looking text is
        # stored as a text value (matching the source data, which keeps
        # these figures as text/shared strings rather than numbers).
        $cell.Value = "'" + $entry[$col]
        # Reset to the default "Normal" style so no extra text-number
        # format is attached to the cell (keeps formatting identical to
        # the rest of the table).
        $cell.Style = "Normal"
    }
}
